$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers need the cell
# number format forced to Text first, otherwise Excel would store
# them as numeric values (losing the zero-padded / dotted text form).
$textForceCells = @("D5", "D9", "D10", "D11", "D19", "D21", "D23", "D24", "D28", "D32", "D40", "D42", "D44", "D45")
foreach ($addr in $textForceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '25.950.10'
$ws.Range('E2').Value = '  +0.00%  '
$ws.Range('D3').Value = '1.638.72'
$ws.Range('E3').Value = '  -0.03%  '
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range('D5').Value = '214.64'
$ws.Range('E6').Value = '  -0.12%  '
$ws.Range('E7').Value = '  +0.04%  '
$ws.Range('E8').Value = '  -0.34%  '
$ws.Range('D9').Value = '0.0637'
$ws.Range('E9').Value = '  -0.31%  '
$ws.Range('D10').Value = '19.51'
$ws.Range('E10').Value = '  -0.71%  '
$ws.Range('D11').Value = '0.0796'
$ws.Range('E11').Value = '  +0.02%  '
$ws.Range('E12').Value = '  -0.12%  '
$ws.Range('D13').Value = '1.619.78'
$ws.Range('E13').Value = '  -1.54%  '
$ws.Range('E14').Value = '  -0.53%  '
$ws.Range('E15').Value = '  +0.84%  '
$ws.Range('E16').Value = '  -0.44%  '
$ws.Range('D17').Value = '25.981.86'
$ws.Range('E17').Value = '  +0.08%  '
$ws.Range('E18').Value = '  +0.02%  '
$ws.Range('D19').Value = '194.15'
$ws.Range('E19').Value = '  +0.02%  '
$ws.Range('E20').Value = '  -0.91%  '
$ws.Range('D21').Value = '9.86'
$ws.Range('E21').Value = '  -0.75%  '
$ws.Range('E22').Value = '  -1.74%  '
$ws.Range('D23').Value = '0.133'
$ws.Range('E23').Value = '  +4.13%  '
$ws.Range('D24').Value = '143.80'
$ws.Range('E24').Value = '  -0.22%  '
$ws.Range('E25').Value = '  -0.04%  '
$ws.Range('E26').Value = '  -0.60%  '
$ws.Range('E27').Value = '  +0.44%  '
$ws.Range('D28').Value = '15.48'
$ws.Range('E28').Value = '  -0.14%  '
$ws.Range('E29').Value = '  -0.16%  '
$ws.Range('E30').Value = '  -1.46%  '
$ws.Range('E31').Value = '  -0.83%  '
$ws.Range('D32').Value = '3.25'
$ws.Range('E32').Value = '  +0.45%  '
$ws.Range('E33').Value = '  -0.79%  '
$ws.Range('E34').Value = '  +0.69%  '
$ws.Range('E35').Value = '  -0.38%  '
$ws.Range('D36').Value = '1.130.51'
$ws.Range('E36').Value = '  -0.74%  '
$ws.Range('E37').Value = '  -1.28%  '
$ws.Range('E38').Value = '  +0.13%  '
$ws.Range('E39').Value = '  -0.71%  '
$ws.Range('D40').Value = '98.53'
$ws.Range('E40').Value = '  -0.85%  '
$ws.Range('E41').Value = '  -0.17%  '
$ws.Range('D42').Value = '0.791'
$ws.Range('E42').Value = '  -1.09%  '
$ws.Range('E43').Value = '  -0.31%  '
$ws.Range('D44').Value = '56.32'
$ws.Range('E44').Value = '  -0.42%  '
$ws.Range('D45').Value = '1.50'
$ws.Range('E45').Value = '  +2.60%  '
$ws.Range('E46').Value = '  -1.57%  '
$ws.Range('E47').Value = '  +0.87%  '
$ws.Range('E48').Value = '  -0.42%  '
$ws.Range('E49').Value = '  -0.18%  '
$ws.Range('E50').Value = '  -1.93%  '
$ws.Range('E51').Value = '  -0.64%  '
